$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.401.75'
$ws.Range('E2').Value = '''  +0.35%  '
$ws.Range('D3').Value = '''1.881.39'
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '''  +0.09%  '
$ws.Range('D5').Value = '''0.7207'
$ws.Range('E5').Value = '''  +1.78%  '
$ws.Range('D6').Value = '''243.47'
$ws.Range('E6').Value = '''  +0.69%  '
$ws.Range('E7').Value = '''  +0.06%  '
$ws.Range('D8').Value = '''0.07967'
$ws.Range('E8').Value = '''  +2.57%  '
$ws.Range('D9').Value = '''0.3154'
$ws.Range('E9').Value = '''  +1.72%  '
$ws.Range('D10').Value = '''24.94'
$ws.Range('E10').Value = '''  -0.13%  '
$ws.Range('D11').Value = '''0.08146'
$ws.Range('E11').Value = '''  -2.85%  '
$ws.Range('D12').Value = '''1.899.05'
$ws.Range('D13').Value = '''94.78'
$ws.Range('E13').Value = '''  +4.02%  '
$ws.Range('D14').Value = '''5.244'
$ws.Range('E14').Value = '''  +0.17%  '
$ws.Range('D15').Value = '''0.7114'
$ws.Range('E15').Value = '''  -0.61%  '
$ws.Range('D16').Value = '''6.405'
$ws.Range('E16').Value = '''  +5.18%  '
$ws.Range('D17').Value = '''0.000008452'
$ws.Range('E17').Value = '''  +2.22%  '
$ws.Range('D18').Value = '''29.415.65'
$ws.Range('E18').Value = '''  +0.35%  '
$ws.Range('D19').Value = '''250.83'
$ws.Range('E19').Value = '''  +4.49%  '
$ws.Range('D20').Value = '''13.32'
$ws.Range('E20').Value = '''  +0.92%  '
$ws.Range('D21').Value = '''2.128.64'
$ws.Range('E21').Value = '''  +0.11%  '
$ws.Range('D22').Value = '''1.002'
$ws.Range('E22').Value = '''  +0.18%  '
$ws.Range('D23').Value = '''7.763'
$ws.Range('E23').Value = '''  +0.15%  '
$ws.Range('D24').Value = '''1.001'
$ws.Range('E24').Value = '''  +0.03%  '
$ws.Range('D25').Value = '''0.1594'
$ws.Range('E25').Value = '''  +0.47%  '
$ws.Range('D26').Value = '''9.086'
$ws.Range('E26').Value = '''  +0.63%  '
$ws.Range('D27').Value = '''162.67'
$ws.Range('E27').Value = '''  +0.21%  '
$ws.Range('D28').Value = '''18.86'
$ws.Range('E28').Value = '''  +1.88%  '
$ws.Range('D29').Value = '''1.507'
$ws.Range('E29').Value = '''  +0.08%  '
$ws.Range('D30').Value = '''4.428'
$ws.Range('E30').Value = '''  +0.54%  '
$ws.Range('D31').Value = '''4.292'
$ws.Range('E31').Value = '''  -0.36%  '
$ws.Range('D32').Value = '''1.224'
$ws.Range('E32').Value = '''  -3.92%  '
$ws.Range('D33').Value = '''0.05328'
$ws.Range('E33').Value = '''  -0.57%  '
$ws.Range('D34').Value = '''1.946'
$ws.Range('E34').Value = '''  +0.43%  '
$ws.Range('D35').Value = '''0.7557'
$ws.Range('E35').Value = '''  +0.61%  '
$ws.Range('E36').Value = '''  +0.45%  '
$ws.Range('D37').Value = '''2.695'
$ws.Range('E37').Value = '''  +0.36%  '
$ws.Range('D38').Value = '''0.01885'
$ws.Range('D39').Value = '''1.275.83'
$ws.Range('E39').Value = '''  +3.04%  '
$ws.Range('D40').Value = '''2.765'
$ws.Range('E40').Value = '''  +1.21%  '
$ws.Range('D41').Value = '''6.471'
$ws.Range('E41').Value = '''  -0.21%  '
$ws.Range('D42').Value = '''113.09'
$ws.Range('E42').Value = '''  +4.08%  '
$ws.Range('D43').Value = '''74.97'
$ws.Range('E43').Value = '''  +3.76%  '
$ws.Range('D44').Value = '''0.9055'
$ws.Range('E44').Value = '''  +1.50%  '
$ws.Range('D45').Value = '''0.00000000130'
$ws.Range('E45').Value = '''  +4.19%  '
$ws.Range('D46').Value = '''1.001'
$ws.Range('E46').Value = '''  +0.02%  '
$ws.Range('D47').Value = '''2.027.56'
$ws.Range('E47').Value = '''  +0.42%  '
$ws.Range('D48').Value = '''1.805'
$ws.Range('E48').Value = '''  +0.72%  '
$ws.Range('D49').Value = '''0.5204'
$ws.Range('E49').Value = '''  +0.07%  '
$ws.Range('D50').Value = '''9.517'
$ws.Range('E50').Value = '''  +0.89%  '
$ws.Range('D51').Value = '''0.4371'
$ws.Range('E51').Value = '''  +0.84%  '
